$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SEGMENT_PULSE sheet (CMO dashboard "SEGMENT PULSE" tab)
# Re-generated zone numbers: market share targets normalized to round
# allocation values (25 / 20), and recomputed demand / awareness-gap /
# high-vs-low-segment summary figures.
# ---------------------------------------------------------------------------
$wsSegment = $wb.Worksheets.Item("SEGMENT_PULSE")

# HIGH SEGMENT ANALYSIS header summary (row 4)
$wsSegment.Range("I4").Value = 45.4
$wsSegment.Range("J4").Value = 12

# HIGH SEGMENT ANALYSIS - Center (row 5)
$wsSegment.Range("B5").Value = 25
$wsSegment.Range("C5").Value = 8800
$wsSegment.Range("D5").Value = -13.67
$wsSegment.Range("F5").Value = 20
$wsSegment.Range("I5").Value = 45.4

# HIGH SEGMENT ANALYSIS - West (row 6)
$wsSegment.Range("B6").Value = 25
$wsSegment.Range("C6").Value = 6000
$wsSegment.Range("D6").Value = -12.42
$wsSegment.Range("F6").Value = 20

# HIGH SEGMENT ANALYSIS - North (row 7)
$wsSegment.Range("B7").Value = 25
$wsSegment.Range("C7").Value = 4900
$wsSegment.Range("D7").Value = -3.739999999999998
$wsSegment.Range("E7").Value = 0
$wsSegment.Range("F7").Value = 20

# LOW SEGMENT ANALYSIS - Center (row 14)
$wsSegment.Range("B14").Value = 25
$wsSegment.Range("C14").Value = 18700
$wsSegment.Range("F14").Value = 20

# LOW SEGMENT ANALYSIS - West (row 15)
$wsSegment.Range("B15").Value = 25
$wsSegment.Range("C15").Value = 12750
$wsSegment.Range("F15").Value = 20

# LOW SEGMENT ANALYSIS - North (row 16)
$wsSegment.Range("B16").Value = 25
$wsSegment.Range("C16").Value = 12600
# D16 now carries a (near-zero, floating point) negative awareness gap, so it
# picks up the same "critical" red highlight style already used by D5:D9.
$wsSegment.Range("D5").Copy()
$wsSegment.Range("D16").PasteSpecial(-4122)
$wsSegment.Range("D16").Value = -3.552713678800501 * 0.000000000000001
$wsSegment.Range("E16").Value = 0
$wsSegment.Range("F16").Value = 20

# HIGH vs LOW SEGMENT GAP summary (rows 27 & 29)
$wsSegment.Range("I27").Value = 20.762
$wsSegment.Range("J27").Value = 20.762
$wsSegment.Range("I29").Value = 12
$wsSegment.Range("J29").Value = 12

# ---------------------------------------------------------------------------
# STRATEGY_COCKPIT sheet (unit economics + zonal allocations)
# TV cost/spot dropped 3787 -> 3000 and Radio cost/spot rose 255 -> 300;
# propagate through the dependent formulas and the zonal "Last Sales" /
# "Avg Comp Price" inputs that were regenerated alongside them.
# ---------------------------------------------------------------------------
$wsStrategy = $wb.Worksheets.Item("STRATEGY_COCKPIT")

# Unit economics cheat sheet (row 3)
$wsStrategy.Range("A3").Value = 3000
$wsStrategy.Range("B3").Value = 300

# Section A: Global allocations - TV spend formula keys off the new TV cost/spot
$wsStrategy.Range("C9").Formula = "=B9*3000.0"

# Section B: Zonal allocations - Center (row 16)
$wsStrategy.Range("B16").Value = 8900
$wsStrategy.Range("H16").Value = 68
$wsStrategy.Range("K16").Formula = "=(C9/5) + (E16*300.0) + (F16*1500) + (MAX(0, F16-5)*1100)"

# Section B: Zonal allocations - West (row 17)
$wsStrategy.Range("B17").Value = 4028
$wsStrategy.Range("H17").Value = 68
$wsStrategy.Range("K17").Formula = "=(C9/5) + (E17*300.0) + (F17*1500) + (MAX(0, F17-5)*1100)"

# Section B: Zonal allocations - North (row 18)
$wsStrategy.Range("B18").Value = 4400
$wsStrategy.Range("H18").Value = 91
$wsStrategy.Range("K18").Formula = "=(C9/5) + (E18*300.0) + (F18*1500) + (MAX(0, F18-5)*1100)"

# Section B: Zonal allocations - East (row 19) - Mkt Cost formula only
$wsStrategy.Range("K19").Formula = "=(C9/5) + (E19*300.0) + (F19*1500) + (MAX(0, F19-5)*1100)"

# Section B: Zonal allocations - South (row 20) - Mkt Cost formula only
$wsStrategy.Range("K20").Formula = "=(C9/5) + (E20*300.0) + (F20*1500) + (MAX(0, F20-5)*1100)"
